$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Broń")
$ws.Range("A48").Value = "Test"
